$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.636.57'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.81%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.860.33'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.01%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.87%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6337'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.54%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.13%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07632'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.79%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3002'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.69'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.09%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07733'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.95%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.930.89'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.53%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6956'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.34%  '

$ws.Range("E14").Value = '  -0.64%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.81'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.07%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001002'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.57%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.102.28'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.27%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.272'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.736.25'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.62%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '235.24'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.54%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.63'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.68%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9995'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.669'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.19%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9999'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.37%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1404'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.73%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.523'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.80'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.71%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.483'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05820'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.259'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.59%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.143'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.57%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.049'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.11%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.909'
$ws.Range("D34").Style = "Normal"

$ws.Range("E35").Value = '  -0.84%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7238'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.23%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.585'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.05%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.255.14'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.89%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.822'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.14%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01811'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.78%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9041'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.42%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.194'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.18%  '

$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '68.99'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.26%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9993'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.12%  '

$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.008.49'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.49%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.49'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.08%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.388'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.76%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.247'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.08%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4074'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.07%  '

$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.724'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.23%  '

$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.00000000117'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.16%  '
